$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("kectest")

# --- Rows whose column-C cluster list was cleared out entirely in the
#     regenerated export (script run produced no cluster membership for
#     these PSTVA rows). ---
$clearRows = @(50, 60, 80, 167, 173, 191, 293, 316, 324, 339)
foreach ($r in $clearRows) {
    $ws.Cells.Item($r, 3).Value = ""
}

# --- Rows whose column-C comma-separated cluster list kept the same
#     members but came out of the regenerated export in a new order. ---
$smpValue = "B126IMD,B126LOG,B126GRG,B126HNA,B126SMP,B126MAS,B126FRC"
$smpRows = @(99, 100, 101, 102, 103, 104, 105, 106, 107, 108, 109, 111, 112, 113, 114, 115, 116, 132, 133, 134, 135, 138)
foreach ($r in $smpRows) {
    $ws.Cells.Item($r, 3).Value = $smpValue
}

$vokValue = "B126AMA,B126BC1,B126GUM,B126CHR,B126PAF,B126GRS,B126PAS,B126VOK,B126SIS,B126CHK,B126JAG,B126MA7,B126ISM"
$vokRows = @(288)
foreach ($r in $vokRows) {
    $ws.Cells.Item($r, 3).Value = $vokValue
}
